# Updated capital structure database
# Refresh the capital-structure metrics (margins, returns, debt, interest
# coverage, etc.) for the Netherlands / "Electronics (Consumer & Office)"
# rows on the "earnings_debt" sheet. The
# historical_growth_net_income_last_5_years column (E) is dropped for both
# rows; historical_growth_revenue_last_5_years (D) and every metric from
# ebitdard_margin (G) through ebit_net_interest_expenses (AQ) are
# refreshed with newly recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Electronics (Consumer & Office) / company "1") ---
# Historical growth columns: net-income-growth column (E) removed; revenue-growth (D) updated
$ws.Range("D2").Value = -0.107
$ws.Range("E2").ClearContents()
$ws.Range("G2").Value = 0.6360450944546009
$ws.Range("H2").Value = 0.05895795246800732
$ws.Range("I2").Value = -0.5135892078670873
$ws.Range("J2").Value = -0.5135892078670873
$ws.Range("K2").Value = -306.1
$ws.Range("L2").Value = -0.4663315051797685
$ws.Range("M2").Value = 19.5
$ws.Range("N2").Value = 0.01450569069404151
$ws.Range("O2").Value = -0.06370467167592289
$ws.Range("R2").Value = -0.0
$ws.Range("S2").Value = 19.5
$ws.Range("U2").Value = 230.7
$ws.Range("V2").Value = 0.1716134791341218
$ws.Range("W2").Value = -0.3819089207735496
$ws.Range("X2").Value = 0.06477607296285352
$ws.Range("Y2").Value = -0.4466849937364031
$ws.Range("Z2").Value = 3.258865635161382
$ws.Range("AA2").Value = -1.673718220107806
$ws.Range("AB2").Value = 0.0630631160655598
$ws.Range("AC2").Value = -1.736781336173366
$ws.Range("AD2").Value = 47.1
$ws.Range("AE2").Value = 0.2197802197802197
$ws.Range("AF2").Value = 47.31978021978022
$ws.Range("AG2").Value = -183.3802197802198
$ws.Range("AH2").Value = 0.03400338288688808
$ws.Range("AI2").Value = 0.08210681264789275
$ws.Range("AJ2").Value = -0.1579611467602895
$ws.Range("AK2").Value = -0.530583694207571
$ws.Range("AL2").Value = 2.7
$ws.Range("AM2").Value = 1.31
$ws.Range("AN2").Value = -1.58714112414072
$ws.Range("AO2").Value = -124.9259259259259
$ws.Range("AP2").Value = 6.179411638368371
$ws.Range("AQ2").Value = -257.4809160305343

# --- Row 3 (TomTom N.V. (ENXTAM:TOM2)) ---
$ws.Range("D3").Value = -0.107
$ws.Range("E3").ClearContents()
$ws.Range("G3").Value = 0.6360450944546009
$ws.Range("H3").Value = 0.05895795246800732
$ws.Range("I3").Value = -0.5135892078670873
$ws.Range("J3").Value = -0.5135892078670873
$ws.Range("K3").Value = -306.1
$ws.Range("L3").Value = -0.4663315051797685
$ws.Range("M3").Value = 19.5
$ws.Range("N3").Value = 0.01450569069404151
$ws.Range("O3").Value = -0.06370467167592289
$ws.Range("R3").Value = 0.0
$ws.Range("S3").Value = 19.5
$ws.Range("U3").Value = 230.7
$ws.Range("V3").Value = 0.1716134791341218
$ws.Range("W3").Value = -0.3819089207735496
$ws.Range("X3").Value = 0.06477607296285352
$ws.Range("Y3").Value = -0.4466849937364031
$ws.Range("Z3").Value = 3.258865635161382
$ws.Range("AA3").Value = -1.673718220107806
$ws.Range("AB3").Value = 0.0630631160655598
$ws.Range("AC3").Value = -1.736781336173366
$ws.Range("AD3").Value = 47.1
$ws.Range("AE3").Value = 0.2197802197802197
$ws.Range("AF3").Value = 47.31978021978022
$ws.Range("AG3").Value = -183.3802197802198
$ws.Range("AH3").Value = 0.03400338288688808
$ws.Range("AI3").Value = 0.08210681264789275
$ws.Range("AJ3").Value = -0.1579611467602895
$ws.Range("AK3").Value = -0.530583694207571
$ws.Range("AL3").Value = 2.7
$ws.Range("AM3").Value = 1.31
$ws.Range("AN3").Value = -1.58714112414072
$ws.Range("AO3").Value = -124.9259259259259
$ws.Range("AP3").Value = 6.179411638368371
$ws.Range("AQ3").Value = -257.4809160305343
